# Update the dSF column (F) with repulled data values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 0
    5  = -7
    7  = -2
    9  = -5
    10 = 0
    11 = -3
    13 = -2
    14 = 7
    15 = 1
    18 = 1
    19 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
